$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell (outside the used range) used to coerce numeric-looking
# strings (e.g. "578.03") into genuine text cells without leaving any
# stray NumberFormat/style behind on the target cells.
$scratch = $ws.Cells.Item(200, 200)

function Set-TextValue {
    param($row, $col, $text)
    $cell = $ws.Cells.Item($row, $col)
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
}

# Row 2
$ws.Cells.Item(2, 4).Value = "61.999.83"
$ws.Cells.Item(2, 5).Value = "  +4.91%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.416.02"
$ws.Cells.Item(3, 5).Value = "  +3.70%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.23%  "

# Row 5
Set-TextValue 5 4 "578.03"
$ws.Cells.Item(5, 5).Value = "  +4.20%  "

# Row 6
Set-TextValue 6 4 "138.59"
$ws.Cells.Item(6, 5).Value = "  +9.29%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.06%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "3.414.48"
$ws.Cells.Item(8, 5).Value = "  +3.46%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  +2.67%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  +2.39%  "

# Row 11
Set-TextValue 11 4 "0.127"
$ws.Cells.Item(11, 5).Value = "  +10.30%  "

# Row 12
Set-TextValue 12 4 "0.394"
$ws.Cells.Item(12, 5).Value = "  +6.99%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "4.001.74"
$ws.Cells.Item(13, 5).Value = "  +3.91%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  +2.17%  "

# Row 15
$ws.Cells.Item(15, 5).Value = "  +8.76%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "3.415.99"
$ws.Cells.Item(16, 5).Value = "  +3.76%  "

# Row 17
Set-TextValue 17 4 "25.48"
$ws.Cells.Item(17, 5).Value = "  +6.56%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "61.998.53"
$ws.Cells.Item(18, 5).Value = "  +4.70%  "

# Row 19
Set-TextValue 19 4 "14.18"
$ws.Cells.Item(19, 5).Value = "  +7.63%  "

# Row 20
Set-TextValue 20 4 "5.89"
$ws.Cells.Item(20, 5).Value = "  +5.12%  "

# Row 21
Set-TextValue 21 4 "9.54"
$ws.Cells.Item(21, 5).Value = "  +7.78%  "

# Row 22
Set-TextValue 22 4 "390.48"
$ws.Cells.Item(22, 5).Value = "  +11.95%  "

# Row 23
Set-TextValue 23 4 "0.572"
$ws.Cells.Item(23, 5).Value = "  +4.05%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "3.552.94"
$ws.Cells.Item(24, 5).Value = "  +3.85%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  +18.73%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  -0.01%  "

# Row 27
Set-TextValue 27 4 "71.51"
$ws.Cells.Item(27, 5).Value = "  +5.06%  "

# Row 28
$ws.Cells.Item(28, 2).Value = "RenderToken"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue 28 4 "7.71"
$ws.Cells.Item(28, 5).Value = "  +6.86%  "

# Row 29
$ws.Cells.Item(29, 2).Value = "Fetch.AI"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue 29 4 "1.59"
$ws.Cells.Item(29, 5).Value = "  +11.04%  "

# Row 30
Set-TextValue 30 4 "1.00"
$ws.Cells.Item(30, 5).Value = "  -0.02%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  +7.59%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  +6.49%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  +5.28%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "3.447.49"
$ws.Cells.Item(34, 5).Value = "  +3.91%  "

# Row 36
Set-TextValue 36 4 "23.59"
$ws.Cells.Item(36, 5).Value = "  +4.62%  "

# Row 37
Set-TextValue 37 4 "5.49"
$ws.Cells.Item(37, 5).Value = "  +4.72%  "

# Row 38
Set-TextValue 38 4 "7.01"
$ws.Cells.Item(38, 5).Value = "  +3.94%  "

# Row 39
Set-TextValue 39 4 "1.56"
$ws.Cells.Item(39, 5).Value = "  +6.24%  "

# Row 40
Set-TextValue 40 4 "162.82"
$ws.Cells.Item(40, 5).Value = "  +2.97%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +6.90%  "

# Row 42
Set-TextValue 42 4 "1.76"
$ws.Cells.Item(42, 5).Value = "  +16.27%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  +7.23%  "

# Row 44
Set-TextValue 44 4 "1.00"
$ws.Cells.Item(44, 5).Value = "  +0.27%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "  +6.25%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  +5.03%  "

# Row 47
Set-TextValue 47 4 "25.22"
$ws.Cells.Item(47, 5).Value = "  +11.62%  "

# Row 48
Set-TextValue 48 4 "41.77"
$ws.Cells.Item(48, 5).Value = "  +3.31%  "

# Row 49
Set-TextValue 49 4 "6.97"
$ws.Cells.Item(49, 5).Value = "  +4.38%  "

# Row 50
Set-TextValue 50 4 "23.11"
$ws.Cells.Item(50, 5).Value = "  +7.68%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "2.382.97"
$ws.Cells.Item(51, 5).Value = "  +12.01%  "

# Clean up the scratch cell so it leaves no structural trace
$scratch.ClearContents()
